$d = $word.ActiveDocument

# Locate the two trailing footer paragraphs that must be removed:
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
# The empty paragraphs immediately before and after them must be preserved.
$startPara = $null
$endPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $startPara = $p
    }
    if ($t -like "*Creative Commons Attribution*") {
        $endPara = $p
    }
}

if (($startPara -ne $null) -and ($endPara -ne $null)) {
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
